$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Gold Answers"
$ws.Range("B31").Value = 10
$ws.Range("E31").Value = "Extractive"
$ws.Range("F31").Value = 34.5

$ws.Range("B32").Value = 10.5
$ws.Range("F32").Value = 39

$ws.Range("B33").Value = 10
$ws.Range("F33").Value = 13

$ws.Range("B34").Value = 14
$ws.Range("F34").Value = 16

$ws.Range("B35").Formula = "=SUM(B31:B34)"
$ws.Range("C35").Formula = "=B35/200"
$ws.Range("F35").Formula = "=SUM(F31:F34)"
$ws.Range("G35").Formula = "=F35/200"

$ws.Range("D38").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
